# corrección bisección 2.0 (lectura de potencias)
# Recompute the bisection iteration table (columns xn, fxn, E) in Sheet1
# with corrected values. Cells are text-typed (as in the original file),
# so a leading apostrophe forces Excel to store them as text rather than
# re-interpreting the numeric-looking strings as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "'0.0"
$ws.Range("C2").Value = "'-11.0"
$ws.Range("B3").Value = "'-2.5"
$ws.Range("C3").Value = "'-3.78125"
$ws.Range("D3").Value = "'2.5"
$ws.Range("B4").Value = "'-3.75"
$ws.Range("C4").Value = "'4.05697572827198"
$ws.Range("D4").Value = "'1.25"
$ws.Range("B5").Value = "'-3.125"
$ws.Range("C5").Value = "'-0.247514006488339"
$ws.Range("D5").Value = "'0.625"
$ws.Range("B6").Value = "'-3.4375"
$ws.Range("C6").Value = "'1.80788665833855"
$ws.Range("D6").Value = "'0.3125"
$ws.Range("B7").Value = "'-3.28125"
$ws.Range("C7").Value = "'0.756021441039556"
$ws.Range("D7").Value = "'0.15625"
$ws.Range("B8").Value = "'-3.203125"
$ws.Range("C8").Value = "'0.248219418534426"
$ws.Range("D8").Value = "'0.078125"
$ws.Range("B9").Value = "'-3.1640625"
$ws.Range("C9").Value = "'-0.0011549191511122"
$ws.Range("D9").Value = "'0.0390625"
$ws.Range("B10").Value = "'-3.18359375"
$ws.Range("C10").Value = "'0.123155220682149"
$ws.Range("D10").Value = "'0.01953125"
$ws.Range("B11").Value = "'-3.173828125"
$ws.Range("C11").Value = "'0.0609059085934192"
$ws.Range("D11").Value = "'0.009765625"
$ws.Range("B12").Value = "'-3.1689453125"
$ws.Range("C12").Value = "'0.0298519360855689"
$ws.Range("D12").Value = "'0.0048828125"
$ws.Range("B13").Value = "'-3.16650390625"
$ws.Range("C13").Value = "'0.0143426190481772"
$ws.Range("D13").Value = "'0.00244140625"
$ws.Range("B14").Value = "'-3.165283203125"
$ws.Range("C14").Value = "'0.0065923776238392"
$ws.Range("D14").Value = "'0.001220703125"
$ws.Range("B15").Value = "'-3.1646728515625"
$ws.Range("C15").Value = "'0.0027183611589549"
$ws.Range("D15").Value = "'0.0006103515625"
$ws.Range("B16").Value = "'-3.16436767578125"
$ws.Range("C16").Value = "'0.0007816289850399"
$ws.Range("D16").Value = "'0.00030517578125"
$ws.Range("B17").Value = "'-3.16421508789062"
$ws.Range("C17").Value = "'-0.0001866680876965"
$ws.Range("D17").Value = "'0.000152587890625"
$ws.Range("B18").Value = "'-3.16429138183594"
$ws.Range("C18").Value = "'0.0002974746974988"
$ws.Range("D18").Value = "'7.62939453125e-05"
$ws.Range("B19").Value = "'-3.16425323486328"
$ws.Range("C19").Value = "'5.54018671081735e-05"
$ws.Range("D19").Value = "'3.814697265625e-05"
$ws.Range("B20").Value = "'-3.16423416137695"
$ws.Range("C20").Value = "'-6.56334697417549e-05"
$ws.Range("D20").Value = "'1.9073486328125e-05"
$ws.Range("B21").Value = "'-3.16424369812012"
$ws.Range("C21").Value = "'-5.11589117913047e-06"
$ws.Range("D21").Value = "'9.5367431640625e-06"
$ws.Range("B22").Value = "'-3.1642484664917"
$ws.Range("C22").Value = "'2.51429654998248e-05"
$ws.Range("D22").Value = "'4.76837158203125e-06"
